$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used by plain (unformatted) data cells, so that when we
# force text formatting (to stop Excel from auto-converting date-looking
# strings like "01-05-2022" into real dates) we can restore the original,
# un-styled look afterwards.
$plainStyle = $ws.Range("A2").Style

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $plainStyle
}

# --- Row 2: AN22-000004 -> AN22-000023 ---
$ws.Range("A2").Value = "AN22-000023"
$ws.Range("B2").Value = "STRAWBERRY MANSION HS - REM 15 SF OF DEBRIS IN CRAWLSPACE 2, REM 20 SF OF GLUE DOTS IN FAN ROOM 215, CAP <1 LF ACPFI IN ATTIC ABOVE AUDITORIUM."
$ws.Range("C2").Value = "Notification"
$ws.Range("E2").Value = "Approved-ASB"
$ws.Range("F2").Value = "3133 RIDGE AVE"
Set-TextValue $ws.Range("G2") "01-05-2022"
Set-TextValue $ws.Range("H2") "01-05-2022"
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
Set-TextValue $ws.Range("K2") "01-07-2022"
Set-TextValue $ws.Range("L2") "01-31-2022"
$ws.Range("N2").ClearContents()
$ws.Range("O2").Value = "Minor Removal Project"
$ws.Range("Q2").Value = "STRAWBERRY MANSION HS"
$ws.Range("X2").Value = 0.5
$ws.Range("Y2").Value = 15
$ws.Range("AB2").Value = 20
$ws.Range("AD2").Value = "No"

# --- Row 3: AN22-000017 -> AN22-000021 ---
$ws.Range("A3").Value = "AN22-000021"
$ws.Range("B3").Value = "Containment for boiler abatement. Removal of pipe insulation using glove bag method. Removal of tranite and wire by butting and wrapping. Removal of fire doors. PCM testing."
$ws.Range("C3").Value = "Notification"
$ws.Range("D3").Value = "Vincent Primavera, III"
$ws.Range("E3").Value = "Approved-ASB"
$ws.Range("F3").Value = "6523-43 LANSDOWNE AVE Lewis C. Cassidy Elementary School"
Set-TextValue $ws.Range("G3") "01-05-2022"
Set-TextValue $ws.Range("H3") "01-05-2022"
$ws.Range("I3").ClearContents()
$ws.Range("J3").ClearContents()
Set-TextValue $ws.Range("K3") "01-12-2022"
Set-TextValue $ws.Range("L3") "02-12-2022"
$ws.Range("N3").ClearContents()
$ws.Range("O3").Value = "Minor Removal Project"
$ws.Range("P3").Value = "Abatement prior to Demolition"
$ws.Range("Q3").Value = "PSD Cassidy Elementary School"
$ws.Range("R3").Value = "School District of Philadelphia"
$ws.Range("S3").Value = "440 N Broad St, 3rd Fl, Rm 3053 Philadelphia, PA 19139"
$ws.Range("T3").Value = "PRIME GROUP REMEDIATION"
$ws.Range("X3").Value = 39
$ws.Range("Y3").Value = 70
$ws.Range("AA3").Value = 132
$ws.Range("AB3").Value = 333
$ws.Range("AD3").Value = "No"

# --- Row 4: AN22-000011 -> AN22-000032 ---
$ws.Range("A4").Value = "AN22-000032"
$ws.Range("B4").Value = "Removal of ACPI and floor tile and mastic"
$ws.Range("C4").Value = "Notification with Alternative Methods"
$ws.Range("D4").Value = "June Huntbach"
$ws.Range("E4").Value = "Online Application Receiv"
$ws.Range("F4").Value = "8060 WILLIAMS AVE"
Set-TextValue $ws.Range("G4") "01-06-2022"
$ws.Range("H4").ClearContents()
Set-TextValue $ws.Range("K4") "01-10-2022"
Set-TextValue $ws.Range("L4") "01-24-2022"
$ws.Range("N4").Value = 2078037
$ws.Range("O4").Value = "Major Removal Project"
$ws.Range("Q4").Value = "FS Edmonds School"
$ws.Range("R4").Value = "The School District of Philadelphia"
$ws.Range("S4").Value = "440 N Broad Street Philadelphia, PA 19130"
$ws.Range("T4").Value = "DIAMOND HUNTBACH, INC."
$ws.Range("X4").Value = 250
$ws.Range("AB4").Value = 100
$ws.Range("AD4").Value = "Yes"

# --- Row 5: AN22-000009 -> AN22-000027 ---
$ws.Range("A5").Value = "AN22-000027"
$ws.Range("B5").Value = "DUNBAR ES - CAP <1 LF ACPI IN ART CLASSROOM 207, 201, CAP <1 LF ACPI & REM 4 SF OF VAT IN SPECIAL EDUCATION 110, CAP <1 LF ACPI IN MAIN HALL H14, BOILER ROOM, BOILER ROOM AT ENTRANCE & BOILER ROOM LEFT CORNER, CAP <1 LF ACPI & ACPFI IN HALL 5 OUTSIDE BOY'S RESTROOM."
$ws.Range("C5").Value = "Notification with Alternative Methods"
$ws.Range("E5").Value = "Online Application Receiv"
$ws.Range("F5").Value = "1750 N 12TH ST"
Set-TextValue $ws.Range("G5") "01-06-2022"
$ws.Range("H5").ClearContents()
Set-TextValue $ws.Range("K5") "01-07-2022"
Set-TextValue $ws.Range("L5") "01-31-2022"
$ws.Range("O5").Value = "Minor Removal Project"
$ws.Range("Q5").Value = "DUNBAR ES"
$ws.Range("X5").Value = 3
$ws.Range("Y5").Value = 0
$ws.Range("AB5").Value = 4

# --- Row 6: delete entirely ---
$ws.Rows(6).Delete()
